$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 2, pushing existing rows down.
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the Waycross/GA entry.
$ws.Range("A2").Value = "Any"
$ws.Range("B2").Value = "Any"
$ws.Range("C2").Value = "Waycross"
$ws.Range("D2").Value = "GA"
$ws.Range("E2").Value = "<50000.0"
$ws.Range("F2").Value = "Skip"

# Update the selection to match the saved state.
$ws.Range("A1:F9").Select() | Out-Null
